$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Unprotect()

# 1. Update the confidentiality / as-of-date notice text in cell A16.
$newline = [char]10
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + $newline + "Model holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."

# 2. Update the Weight (D) and Percent Change (E) values for rows 2-13.
$ws.Range("D2").Value = 0.03102874635442196
$ws.Range("E2").Value = 0.01160805735745973

$ws.Range("D3").Value = 0.0238936380363264
$ws.Range("E3").Value = -0.001404823226410779

$ws.Range("D4").Value = 0.05247347009513081
$ws.Range("E4").Value = 0.002545118000925495

$ws.Range("D5").Value = 0.1392088866521911
$ws.Range("E5").Value = 0.001292824822236582

$ws.Range("D6").Value = 0.0315180769011938
$ws.Range("E6").Value = 0.008474576271186196

$ws.Range("D7").Value = 0.1160815609632539
$ws.Range("E7").Value = 0.008985544992837768

$ws.Range("D8").Value = 0.1021014193799655
$ws.Range("E8").Value = 0.005894271504881088

$ws.Range("D9").Value = 0.02942375168337564
$ws.Range("E9").Value = 0.008046214153084597

$ws.Range("D10").Value = 0.1275520831721475
$ws.Range("E10").Value = 0.003359462486002274

$ws.Range("D11").Value = 0.2442939011392796
$ws.Range("E11").Value = 0.0177744862062581

$ws.Range("D12").Value = 0.1024244656227139
$ws.Range("E12").Value = 0.01893791487700125

$ws.Range("E13").Value = 0.00949927356255742

# Restore sheet protection (it was removed above so the locked cells
# could be edited). The original file does not show protection being
# removed in the target diff, so we turn it back on here.
$ws.Protect()
